$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: push the old 2022-Q3 summary row down to row 3
#    and write the new 2022-Q4 summary numbers into row 2.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Duplicate row 2 (value + style) down into row 3 so the old quarter's totals
# are preserved as history underneath the newest quarter.
$summary.Range("A2:D2").Copy($summary.Range("A3:D3"))
$summary.Range("A3").Value = 1

# Row 2 now becomes the latest quarter (2022-Q4) entry.
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("D2").Value = 0.18

# ---------------------------------------------------------------------------
# 2) Archive the existing "2022-Q3" sheet as-is (values + formatting) into a
#    new sheet placed right after it, so the quarterly detail history is kept.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)
$q3.Copy($null, $q3)

# ---------------------------------------------------------------------------
# 3) Turn the original sheet into the new "2022-Q4" detail sheet: rename it
#    and replace its contents with the Q4 fund holdings.
# ---------------------------------------------------------------------------
$q3.Name = "2022-Q4"
$q3.Cells.Clear()

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "001735"
$q3.Range("C2").Value = "广发百发大数据策略成长灵活配置混合E"
$q3.Range("D2").Value = "8.88"
$q3.Range("E2").Value = "90.16"
$q3.Range("F2").Value = "1.36"
$q3.Range("G2").Value = "0.1208"
$q3.Range("H2").Value = 10

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "001734"
$q3.Range("C3").Value = "广发百发大数据策略成长灵活配置混合A"
$q3.Range("D3").Value = "4.29"
$q3.Range("E3").Value = "90.16"
$q3.Range("F3").Value = "1.36"
$q3.Range("G3").Value = "0.0583"
$q3.Range("H3").Value = 10

# Match the header/id-column styling used elsewhere in the workbook (bold,
# bordered, center/top aligned) by copying the format from the "总计" header.
$summary.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$q3.Range("A2").PasteSpecial(-4122)
$q3.Range("A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) The archived copy currently sits at the end named "2022-Q3 (2)"; restore
#    its original name now that the source sheet has been renamed away.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(3).Name = "2022-Q3"

$q3.Range("A1").Select()

Write-Output "done"
